{"js": "// The document body consists of a single paragraph whose single run\n// contains a long sequence of text segments (\"<w:t>\") separated by\n// manual line breaks (\"<w:br/>\"). The edit rewrites almost all of the\n// text (new title line, reworded legal clauses, renamed parties, etc.)\n// while preserving that same text/line-break structure.\n//\n// Word represents a manual line break as a vertical-tab character\n// (U+000B, \"\\v\") inside Range/Body.text. Inserting one string that mixes\n// literal text with \"\\v\" characters therefore reproduces the exact\n// \"<w:t>...</w:t><w:br/>\" run layout the diff shows.\nconst body = context.document.body;\n\n// Wipe the existing paragraph content \u2026\nbody.clear();\nawait context.sync();\n\n// \u2026 and insert the full replacement text (with \"\\v\" standing in for\n// each \"<w:br/>\") in one shot so it lands in a single run, matching the\n// target markup.\nconst newText = \"Erteilung einer Erlaubnis f\u00fcr eine Musikkneipe in den R\u00e4umen der ehemaligen Gastst\u00e4tte \\\"Walfisch\\\" in Kehl\\u000b\\u000bRechtsgrundlage\\u000bDie Rechtsgrundlage f\u00fcr die Erteilung der Erlaubnis k\u00f6nnte \u00a7 12 GastG sein.\\u000b\\u000bMaterielle Voraussetzungen\\u000b\\u000bTatbestandsvoraussetzung\\u000b\\u000bErlaubnispflicht\\u000bDie Erlaubnispflicht ergibt sich aus \u00a7 2 GastG.\\u000b\\u000bZuverl\u00e4ssigkeit\\u000bDie Zuverl\u00e4ssigkeit des Antragstellers ist gem\u00e4\u00df \u00a7 4 GastG zu pr\u00fcfen.\\u000b\\u000bSachkunde\\u000bDie Sachkunde des Antragstellers ist gem\u00e4\u00df \u00a7 5 GastG zu pr\u00fcfen.\\u000b\\u000bRechtsfolgenseite\\u000b\\u000bDer Pflichtige\\u000bAls Pflichtiger kommt Michael Graeter in Betracht, da er die Erlaubnis beantragt hat.\\u000b\\u000bErmessen\\u000bDie Stadt Kehl hat gem\u00e4\u00df \u00a7 12 GastG ein Ermessen, das nach \u00a7 40 LVwVfG ausge\u00fcbt wird. Die Erlaubniserteilung ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, wenn die materiellen Voraussetzungen erf\u00fcllt sind und die \u00f6ffentlichen Interessen nicht entgegenstehen.\\u000b\\u000bUnm\u00f6glichkeit\\u000bEs k\u00f6nnte eine Unm\u00f6glichkeit vorliegen, wenn die Herrentoilette nicht den Anforderungen entspricht. In diesem Fall k\u00f6nnte die Erlaubnis nur unter der Bedingung erteilt werden, dass die Toilette entsprechend umgebaut wird.\\u000b\\u000bBestimmtheit\\u000bNach \u00a7 37 LVwVfG muss die Erlaubnis bestimmt genug formuliert werden.\\u000b\\u000bFormelle Voraussetzungen\\u000b\\u000bZust\u00e4ndigkeit\\u000b\\u000bSachliche Zust\u00e4ndigkeit\\u000bNach \u00a7 2 Abs. 1 GastG ist die Stadt Kehl sachlich zust\u00e4ndig.\\u000b\\u000b\u00d6rtliche Zust\u00e4ndigkeit\\u000b\u00d6rtlich zust\u00e4ndig ist die Stadt Kehl gem\u00e4\u00df \u00a7 3 Abs.1 Nr.1 LVwVfG.\\u000b\\u000bVerfahren\\u000b\\u000bBeteiligte\\u000bNach \u00a7 28 Abs. 1 LVwVfG ist Michael Graeter beteiligt.\\u000b\\u000bAusgeschlossene Personen/Befangenheit\\u000bEs liegen keine Anhaltspunkte f\u00fcr eine Befangenheit vor.\\u000b\\u000bBeteiligung anderer Beh\u00f6rden\\u000bEs sind keine anderen Beh\u00f6rden zu beteiligen.\\u000b\\u000bAnh\u00f6rung\\u000bNach \u00a7 28 Abs. 1 LVwVfG ist Michael Graeter die Gelegenheit zur \u00c4u\u00dferung zu geben.\\u000b\\u000bForm\\u000b\\u000bFormwahl\\u000bNach \u00a7 37 Abs. 2 LVwVfG kann die Erlaubnis hier schriftlich erfolgen.\\u000b\\u000bBegr\u00fcndungspflicht\\u000bNach \u00a7 39 Abs. 1 LVwVfG ist die schriftliche Erlaubnis auch schriftlich zu begr\u00fcnden.\\u000b\\u000bRechtsbehelfsbelehrung\\u000bNach \u00a7 37 Abs. 6 LVwVfG ist eine Rechtsbehelfsbelehrung beizuf\u00fcgen.\\u000b\\u000bBekanntgabe\\u000bNach \u00a7 43 Abs. 1 LVwVfG wird ein Verwaltungsakt durch Bekanntgabe wirksam.\\u000b\\u000bDem Michael Graeter sollte die Erlaubnis mittels PZU nach \u00a7 3 LVwZG zugestellt und damit bekanntgegeben werden. Es sollte darauf hingewiesen werden, dass die Erlaubnis unter der Bedingung erteilt wird, dass die Herrentoilette entsprechend umgebaut wird. Au\u00dferdem sollte darauf hingewiesen werden, dass bei L\u00e4rmbel\u00e4stigungen in der Umgebung Ma\u00dfnahmen ergriffen werden m\u00fcssen.\";\nbody.insertText(newText, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# The document body is a single paragraph whose single run contains a\n# long run of text segments separated by manual line breaks (<w:br/>).\n# The edit rewrites nearly all of that text (new title line, reworded\n# legal clauses, renamed parties, etc.) while keeping the same\n# text/line-break layout, so the whole story is replaced in one go.\n#\n# Word represents a manual line break as a vertical-tab character\n# (charcode 11, written \"`v\" in a PowerShell double-quoted string) inside\n# Range.Text. Assigning one string that mixes literal text with \"`v\"\n# therefore reproduces the exact \"<w:t>...</w:t><w:br/>\" run layout shown\n# in the diff.\n$d = $word.ActiveDocument\n\n$d.Content.Text = \"Erteilung einer Erlaubnis f\u00fcr eine Musikkneipe in den R\u00e4umen der ehemaligen Gastst\u00e4tte `\"Walfisch`\" in Kehl`v`vRechtsgrundlage`vDie Rechtsgrundlage f\u00fcr die Erteilung der Erlaubnis k\u00f6nnte \u00a7 12 GastG sein.`v`vMaterielle Voraussetzungen`v`vTatbestandsvoraussetzung`v`vErlaubnispflicht`vDie Erlaubnispflicht ergibt sich aus \u00a7 2 GastG.`v`vZuverl\u00e4ssigkeit`vDie Zuverl\u00e4ssigkeit des Antragstellers ist gem\u00e4\u00df \u00a7 4 GastG zu pr\u00fcfen.`v`vSachkunde`vDie Sachkunde des Antragstellers ist gem\u00e4\u00df \u00a7 5 GastG zu pr\u00fcfen.`v`vRechtsfolgenseite`v`vDer Pflichtige`vAls Pflichtiger kommt Michael Graeter in Betracht, da er die Erlaubnis beantragt hat.`v`vErmessen`vDie Stadt Kehl hat gem\u00e4\u00df \u00a7 12 GastG ein Ermessen, das nach \u00a7 40 LVwVfG ausge\u00fcbt wird. Die Erlaubniserteilung ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, wenn die materiellen Voraussetzungen erf\u00fcllt sind und die \u00f6ffentlichen Interessen nicht entgegenstehen.`v`vUnm\u00f6glichkeit`vEs k\u00f6nnte eine Unm\u00f6glichkeit vorliegen, wenn die Herrentoilette nicht den Anforderungen entspricht. In diesem Fall k\u00f6nnte die Erlaubnis nur unter der Bedingung erteilt werden, dass die Toilette entsprechend umgebaut wird.`v`vBestimmtheit`vNach \u00a7 37 LVwVfG muss die Erlaubnis bestimmt genug formuliert werden.`v`vFormelle Voraussetzungen`v`vZust\u00e4ndigkeit`v`vSachliche Zust\u00e4ndigkeit`vNach \u00a7 2 Abs. 1 GastG ist die Stadt Kehl sachlich zust\u00e4ndig.`v`v\u00d6rtliche Zust\u00e4ndigkeit`v\u00d6rtlich zust\u00e4ndig ist die Stadt Kehl gem\u00e4\u00df \u00a7 3 Abs.1 Nr.1 LVwVfG.`v`vVerfahren`v`vBeteiligte`vNach \u00a7 28 Abs. 1 LVwVfG ist Michael Graeter beteiligt.`v`vAusgeschlossene Personen/Befangenheit`vEs liegen keine Anhaltspunkte f\u00fcr eine Befangenheit vor.`v`vBeteiligung anderer Beh\u00f6rden`vEs sind keine anderen Beh\u00f6rden zu beteiligen.`v`vAnh\u00f6rung`vNach \u00a7 28 Abs. 1 LVwVfG ist Michael Graeter die Gelegenheit zur \u00c4u\u00dferung zu geben.`v`vForm`v`vFormwahl`vNach \u00a7 37 Abs. 2 LVwVfG kann die Erlaubnis hier schriftlich erfolgen.`v`vBegr\u00fcndungspflicht`vNach \u00a7 39 Abs. 1 LVwVfG ist die schriftliche Erlaubnis auch schriftlich zu begr\u00fcnden.`v`vRechtsbehelfsbelehrung`vNach \u00a7 37 Abs. 6 LVwVfG ist eine Rechtsbehelfsbelehrung beizuf\u00fcgen.`v`vBekanntgabe`vNach \u00a7 43 Abs. 1 LVwVfG wird ein Verwaltungsakt durch Bekanntgabe wirksam.`v`vDem Michael Graeter sollte die Erlaubnis mittels PZU nach \u00a7 3 LVwZG zugestellt und damit bekanntgegeben werden. Es sollte darauf hingewiesen werden, dass die Erlaubnis unter der Bedingung erteilt wird, dass die Herrentoilette entsprechend umgebaut wird. Au\u00dferdem sollte darauf hingewiesen werden, dass bei L\u00e4rmbel\u00e4stigungen in der Umgebung Ma\u00dfnahmen ergriffen werden m\u00fcssen.\"\n"}
